$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "For B5" adjustment-total block (rows 36-41), mirroring the existing
# "For B2" block in rows 27-32.
#
# NOTE: shared-string insertion order matters for matching the target
# sharedStrings.xml table (new uniques get appended in first-use order), so
# A40's "bounds-adjustment" is written before A36's "For B5: -6 6" label.
$ws.Range("A40").Value = "bounds-adjustment"
$ws.Range("A36").Value = "For B5: -6 6"

$ws.Range("A37").Value = "center of mass adjustment"
$ws.Range("C37").Formula = "=0.01*SUM(B2:B5)/4"
$ws.Range("D37").Formula = "=0.01*SUM(C2:C5)/4"

$ws.Range("A38").Value = "velocity-adjustment"
$ws.Range("C38").Formula = "=SUM(D2:D5)/4"
$ws.Range("D38").Formula = "=SUM(E2:E5)/4"

$ws.Range("A39").Value = "avoidance-adjustment"
$ws.Range("C39").Formula = "=-6.5-B5"
$ws.Range("D39").Formula = "=6-C5"

$ws.Range("C40").Value = 1
$ws.Range("D40").Value = 0

$ws.Range("C41").Formula = "=SUM(C37:C40)"
$ws.Range("D41").Formula = "=SUM(D37:D40)"

# Match the saved selection/active cell from the edit.
$ws.Range("C40").Select()
